$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "national_id" column (M), mirroring the "status" column (L) ---
#
# Copy the formatting (fill/border/number-format) from the whole status
# column (L1:L10) onto the brand-new column M first. Doing the format
# copy *before* writing the values means the "Text" number format that
# already lives on the status column carries over, so the numeric-looking
# values we set next ("123") are stored as text/shared-strings exactly
# like the rest of the sheet, instead of being auto-detected as numbers.
$ws.Range("L1:L10").Copy() | Out-Null
$ws.Range("M1:M10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header
$ws.Range("M1").Value = "national_id"

# Data rows
$ws.Range("M2").Value = "123"
# M3 is intentionally left blank (only the copied formatting applies)

# --- Column widths ---
# Columns D:E grew slightly (23.4531 -> 23.5 characters)
$ws.Range("D1:E1").EntireColumn.ColumnWidth = 22.59

# The new column M takes on the same width as columns K:L (~36.6719 characters)
$ws.Range("M1").EntireColumn.ColumnWidth = 35.751
